# ----------------------------------------------------------------------------
# Refresh the cryptocurrency price / 1h-volume snapshot (and a handful of
# re-ranked coin rows) to match the latest GitHub Actions scrape.
#
# Commit: "Updated cryptos list on Thu Sep 21 09:28:34 UTC 2023 with GitHub
# Actions"
#
# Price/Volume cells are stored as TEXT in this workbook (e.g. "19.83",
# "0.0625", "26.956.01" thousand-grouped, or "  -1.30%  " padded with
# whitespace). Excel auto-coerces plain numeric-looking strings assigned via
# .Value into real numbers, which would silently change the stored cell type
# and drop the original text formatting (leading/trailing spaces, the
# dotted-thousands grouping, etc). For every "Price" (column D) cell we
# therefore force the destination NumberFormat to Text ("@") before the
# write, then clear the format straight back off the cell afterwards so no
# extra per-cell style sticks around. "Volume(1h)" (column E) values never
# look like bare numbers (they carry padding/sign/percent), so they're
# assigned directly. Coin name (B) / link (C) cells are plain text too.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '26.932.61'
$c.ClearFormats()
$ws.Range('E2').Value = '  -0.86%  '

# Row 3
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.616.98'
$c.ClearFormats()
$ws.Range('E3').Value = '  -1.31%  '

# Row 4
$ws.Range('E4').Value = '  +0.03%  '

# Row 5
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '212.38'
$c.ClearFormats()
$ws.Range('E5').Value = '  -2.19%  '

# Row 6
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.510'
$c.ClearFormats()
$ws.Range('E6').Value = '  -1.39%  '

# Row 7
$ws.Range('E7').Value = '  +0.09%  '

# Row 8
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.249'
$c.ClearFormats()
$ws.Range('E8').Value = '  -1.84%  '

# Row 9
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.0623'
$c.ClearFormats()
$ws.Range('E9').Value = '  -0.37%  '

# Row 10
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '19.83'
$c.ClearFormats()
$ws.Range('E10').Value = '  -1.46%  '

# Row 11
$ws.Range('E11').Value = '  -1.39%  '

# Row 12
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '1.846.51'
$c.ClearFormats()
$ws.Range('E12').Value = '  -1.13%  '

# Row 13
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '1.620.95'
$c.ClearFormats()
$ws.Range('E13').Value = '  -0.93%  '

# Row 14
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '4.08'
$c.ClearFormats()
$ws.Range('E14').Value = '  -1.38%  '

# Row 15
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.534'
$c.ClearFormats()
$ws.Range('E15').Value = '  -2.02%  '

# Row 16
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '26.921.57'
$c.ClearFormats()
$ws.Range('E16').Value = '  -0.90%  '

# Row 17
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '63.89'
$c.ClearFormats()
$ws.Range('E17').Value = '  -3.69%  '

# Row 18
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '0.0₃0731'
$c.ClearFormats()
$ws.Range('E18').Value = '  -1.04%  '

# Row 19
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '211.26'
$c.ClearFormats()
$ws.Range('E19').Value = '  -2.56%  '

# Row 20
$ws.Range('E20').Value = '  -0.05%  '

# Row 21
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '6.76'
$c.ClearFormats()
$ws.Range('E21').Value = '  -1.63%  '

# Row 22
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.31'
$c.ClearFormats()
$ws.Range('E22').Value = '  -2.71%  '

# Row 23
$ws.Range('E23').Value = '  -7.65%  '

# Row 24
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '8.93'
$c.ClearFormats()
$ws.Range('E24').Value = '  -2.46%  '

# Row 25
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '146.48'
$c.ClearFormats()
$ws.Range('E25').Value = '  -0.87%  '

# Row 26
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '7.47'
$c.ClearFormats()
$ws.Range('E26').Value = '  +0.99%  '

# Row 27
$ws.Range('E27').Value = '  +0.01%  '

# Row 28
$ws.Range('E28').Value = '  -4.31%  '

# Row 29
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '15.40'
$c.ClearFormats()
$ws.Range('E29').Value = '  -1.84%  '

# Row 30
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.0504'
$c.ClearFormats()
$ws.Range('E30').Value = '  -0.99%  '

# Row 31
$ws.Range('E31').Value = '  -1.80%  '

# Row 32
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '3.27'
$c.ClearFormats()
$ws.Range('E32').Value = '  -3.32%  '

# Row 33
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.702'
$c.ClearFormats()
$ws.Range('E33').Value = '  +27.82%  '

# Row 34
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '2.97'
$c.ClearFormats()
$ws.Range('E34').Value = '  -1.90%  '

# Row 35
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.324.80'
$c.ClearFormats()
$ws.Range('E35').Value = '  +1.56%  '

# Row 36
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '1.54'
$c.ClearFormats()
$ws.Range('E36').Value = '  -1.84%  '

# Row 37
$ws.Range('E37').Value = '  -0.62%  '

# Row 38
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.0174'
$c.ClearFormats()
$ws.Range('E38').Value = '  -1.51%  '

# Row 39
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.831'
$c.ClearFormats()
$ws.Range('E39').Value = '  -2.77%  '

# Row 40
$ws.Range('E40').Value = '  +0.01%  '

# Row 41 (MXToken -> TrustWalletToken)
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.793'
$c.ClearFormats()
$ws.Range('E41').Value = '  -2.16%  '

# Row 42 (TrustWalletToken -> MXToken)
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '2.21'
$c.ClearFormats()
$ws.Range('E42').Value = '  -2.85%  '

# Row 43
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '5.30'
$c.ClearFormats()
$ws.Range('E43').Value = '  -0.87%  '

# Row 44
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '63.49'
$c.ClearFormats()
$ws.Range('E44').Value = '  +1.39%  '

# Row 45
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '1.756.46'
$c.ClearFormats()
$ws.Range('E45').Value = '  -1.23%  '

# Row 46
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '89.52'
$c.ClearFormats()
$ws.Range('E46').Value = '  -1.48%  '

# Row 47
$ws.Range('E47').Value = '  +0.60%  '

# Row 48
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.812'
$c.ClearFormats()
$ws.Range('E48').Value = '  +8.33%  '

# Row 49 (Cronos -> BabyDogeCoin)
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.0₆0103'
$c.ClearFormats()
$ws.Range('E49').Value = '  -3.98%  '

# Row 50 (Algorand -> Cronos)
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.0515'
$c.ClearFormats()
$ws.Range('E50').Value = '  -0.04%  '

# Row 51 (EnergySwap -> Algorand)
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.0982'
$c.ClearFormats()
$ws.Range('E51').Value = '  +2.29%  '
